# Apply the "31 Nov 2021 2nd commit" edit to the CCPAUrls sheet:
#  1. Flip the "Status" column (C) from OFF -> ON for rows 5 through 91.
#  2. Update the saved view state (scroll/selection) of the sheet so that
#     the frozen pane keeps its top-left cell at A2 and the active
#     selection moves to H91.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CCPAUrls")

# Make sure we're working on/looking at the right sheet before touching
# the selection, exactly like a user would in the Excel UI.
$ws.Activate() | Out-Null

# --- 1. Update column C (Status) values for rows 5-91 -----------------
for ($row = 5; $row -le 91; $row++) {
    $ws.Cells.Item($row, 3).Value = "ON"
}

# --- 2. Update the view / selection state ------------------------------
# Re-select the frozen top-left anchor first so the pane's recorded
# top-left cell lines up with row 2 (just beneath the frozen header row),
# then move the active selection/cell to H91.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("H91").Select() | Out-Null
